# Bento "startup" tab: the multi-filter scripts were updated for the new
# data set availability - the grouped_recurrence_score filter value changes
# from "51-100" to "0-5" everywhere it is used (the per-tab Neo4j query in
# column B, and the shared StatQuery count query in column C, for the
# CasesTab/SamplesTab/FilesTab rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

$cellsToUpdate = @("B2", "C2", "B3", "C3", "B4", "C4")
foreach ($cellRef in $cellsToUpdate) {
    $cell = $ws.Range($cellRef)
    $queryText = [string]$cell.Value()
    $updatedText = $queryText.Replace('"51-100"', '"0-5"')
    $cell.Value = $updatedText
}

# Scroll the view back up and select C2, matching where the author landed
# after editing the CasesTab row.
$excel.Goto($ws.Range("C2"), $true)
